$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("B2").Value = 1.256650447503432
$ws.Range("D2").Value = 0.2171995174122117
$ws.Range("E2").Value = 0.2504361705155205
$ws.Range("F2").Value = 1.297979492647961
$ws.Range("G2").Value = 0.002445197301900778
$ws.Range("I2").Value = 0.9601397729862136
$ws.Range("J2").Value = 0.3512477623970511
$ws.Range("L2").Value = 0.5328880917857077
$ws.Range("M2").Value = 0.4077681410099032
$ws.Range("O2").Value = 3.175269015285494
# Row 3
$ws.Range("B3").Value = 1.162248352745678
$ws.Range("D3").Value = 0.2164623856048351
$ws.Range("E3").Value = 0.2445548809568407
$ws.Range("F3").Value = 1.310031483620634
$ws.Range("G3").Value = 0.002448498685657444
$ws.Range("I3").Value = 0.9765614704292247
$ws.Range("J3").Value = 0.337761140247693
$ws.Range("L3").Value = 0.4853575348940069
$ws.Range("M3").Value = 0.3745916773015168
$ws.Range("O3").Value = 3.184452731511499
# Row 4
$ws.Range("B4").Value = 1.104219100402702
$ws.Range("D4").Value = 0.2160452416684251
$ws.Range("E4").Value = 0.240975147994952
$ws.Range("F4").Value = 1.318381942020856
$ws.Range("G4").Value = 0.002450635166373738
$ws.Range("I4").Value = 0.9873384794979252
$ws.Range("J4").Value = 0.3295105323216632
$ws.Range("L4").Value = 0.4561198320986364
$ws.Range("M4").Value = 0.3541913984188128
$ws.Range("O4").Value = 3.192565907845903
# Row 5
$ws.Range("B5").Value = 1.08055676267702
$ws.Range("D5").Value = 0.2158842414296842
$ws.Range("E5").Value = 0.239524473450814
$ws.Range("F5").Value = 1.32202350526952
$ws.Range("G5").Value = 0.00245153339765063
$ws.Range("I5").Value = 0.991904650770941
$ws.Range("J5").Value = 0.3261563297484145
$ws.Range("L5").Value = 0.4441925401876858
$ws.Range("M5").Value = 0.3458712077777477
$ws.Range("O5").Value = 3.19649315448126
# Row 6
$ws.Range("B6").Value = 1.076626795647996
$ws.Range("D6").Value = 0.2158580525516882
$ws.Range("E6").Value = 0.2392840848989444
$ws.Range("F6").Value = 1.322642590468803
$ws.Range("G6").Value = 0.002451684217478883
$ws.Range("I6").Value = 0.9926733947253155
$ws.Range("J6").Value = 0.3255998607587429
$ws.Range("L6").Value = 0.4422112804160747
$ws.Range("M6").Value = 0.3444892462595064
$ws.Range("O6").Value = 3.197182748234866
# Row 7
$ws.Range("B7").Value = 1.103900040163722
$ws.Range("D7").Value = 0.2160430338674288
$ws.Range("E7").Value = 0.2409555506459711
$ws.Range("F7").Value = 1.318430087429114
$ws.Range("G7").Value = 0.002450647168194722
$ws.Range("I7").Value = 0.9873993542130997
$ws.Range("J7").Value = 0.3294652635199071
$ws.Range("L7").Value = 0.4559590267529643
$ws.Range("M7").Value = 0.3540792166275608
$ws.Range("O7").Value = 3.192616358924255
# Row 8
$ws.Range("B8").Value = 1.224115252644822
$ws.Range("D8").Value = 0.2169380383392721
$ws.Range("E8").Value = 0.24840190124781
$ws.Range("F8").Value = 1.30193755015253
$ws.Range("G8").Value = 0.002446312960007197
$ws.Range("I8").Value = 0.9656578916760097
$ws.Range("J8").Value = 0.3465915256707035
$ws.Range("L8").Value = 0.5165111890666196
$ws.Range("M8").Value = 0.3963354169588342
$ws.Range("O8").Value = 3.177921324384471
# Row 9
$ws.Range("B9").Value = 1.459273310236938
$ws.Range("D9").Value = 0.2189714544029684
$ws.Range("E9").Value = 0.2632454320157862
$ws.Range("F9").Value = 1.277152559698152
$ws.Range("G9").Value = 0.002438677935320098
$ws.Range("I9").Value = 0.9285324339611414
$ws.Range("J9").Value = 0.3804008965574752
$ws.Range("L9").Value = 0.6347971260456688
$ws.Range("M9").Value = 0.4789420104503535
$ws.Range("O9").Value = 3.168792056765938
# Row 10
$ws.Range("B10").Value = 1.631625493301499
$ws.Range("D10").Value = 0.2206312724491042
$ws.Range("E10").Value = 0.2742886722565245
$ws.Range("F10").Value = 1.263571573650943
$ws.Range("G10").Value = 0.002433590036408743
$ws.Range("I10").Value = 0.9046186492261015
$ws.Range("J10").Value = 0.4053601775859761
$ws.Range("L10").Value = 0.7213917171540061
$ws.Range("M10").Value = 0.5394536667439382
$ws.Range("O10").Value = 3.174168498705797
# Row 11
$ws.Range("B11").Value = 1.70993023340219
$ws.Range("D11").Value = 0.2214216013666359
$ws.Range("E11").Value = 0.2793405415612753
$ws.Range("F11").Value = 1.258402994928403
$ws.Range("G11").Value = 0.002431387525838223
$ws.Range("I11").Value = 0.8944705854440151
$ws.Range("J11").Value = 0.4167372645420073
$ws.Range("L11").Value = 0.7607124624745154
$ws.Range("M11").Value = 0.5669387801064403
$ws.Range("O11").Value = 3.179256128727957
# Row 12
$ws.Range("B12").Value = 1.739566596871043
$ws.Range("D12").Value = 0.2217258808594877
$ws.Range("E12").Value = 0.2812574343280048
$ws.Range("F12").Value = 1.256591325854032
$ws.Range("G12").Value = 0.002430569511320345
$ws.Range("I12").Value = 0.8907329131093711
$ws.Range("J12").Value = 0.4210484345163934
$ws.Range("L12").Value = 0.7755912287476292
$ws.Range("M12").Value = 0.5773401652206189
$ws.Range("O12").Value = 3.181563920981887
# Row 13
$ws.Range("B13").Value = 1.733184611401953
$ws.Range("D13").Value = 0.2216601276412646
$ws.Range("E13").Value = 0.2808444287156888
$ws.Range("F13").Value = 1.256975021088891
$ws.Range("G13").Value = 0.002430744973740316
$ws.Range("I13").Value = 0.8915332075254021
$ws.Range("J13").Value = 0.4201198235063401
$ws.Range("L13").Value = 0.7723873294804946
$ws.Range("M13").Value = 0.5751003452108279
$ws.Range("O13").Value = 3.181049920893344
# Row 14
$ws.Range("B14").Value = 1.71236876395443
$ws.Range("D14").Value = 0.2214465348476296
$ws.Range("E14").Value = 0.2794981692825189
$ws.Range("F14").Value = 1.258251028321595
$ws.Range("G14").Value = 0.002431319906682981
$ws.Range("I14").Value = 0.8941609761091378
$ws.Range("J14").Value = 0.4170918909141079
$ws.Range("L14").Value = 0.7619367757902751
$ws.Range("M14").Value = 0.5677946446206192
$ws.Range("O14").Value = 3.179438344472231
# Row 15
$ws.Range("B15").Value = 1.699616330499794
$ws.Range("D15").Value = 0.22131635197848
$ws.Range("E15").Value = 0.2786740431680741
$ws.Range("F15").Value = 1.259051587516723
$ws.Range("G15").Value = 0.002431674154668861
$ws.Range("I15").Value = 0.8957842619967771
$ws.Range("J15").Value = 0.4152375627305105
$ws.Range("L15").Value = 0.7555340330640092
$ws.Range("M15").Value = 0.5633188108589025
$ws.Range("O15").Value = 3.178500892252856
# Row 16
$ws.Range("B16").Value = 1.626505954879235
$ws.Range("D16").Value = 0.2205803260856598
$ws.Range("E16").Value = 0.2739590717494593
$ws.Range("F16").Value = 1.263929699754719
$ws.Range("G16").Value = 0.002433736222645613
$ws.Range("I16").Value = 0.9052965594687663
$ws.Range("J16").Value = 0.4046170871805543
$ws.Range("L16").Value = 0.7188205005178645
$ws.Range("M16").Value = 0.5376565545475529
$ws.Range("O16").Value = 3.173889284093065
# Row 17
$ws.Range("B17").Value = 1.581628510129917
$ws.Range("D17").Value = 0.2201377795250323
$ws.Range("E17").Value = 0.2710736791445925
$ws.Range("F17").Value = 1.267181098916012
$ws.Range("G17").Value = 0.002435029865546064
$ws.Range("I17").Value = 0.9113192288883525
$ws.Range("J17").Value = 0.3981073767827894
$ws.Range("L17").Value = 0.6962789956161259
$ws.Range("M17").Value = 0.5219024202266667
$ws.Range("O17").Value = 3.17173779658043
# Row 18
$ws.Range("B18").Value = 1.555806959528582
$ws.Range("D18").Value = 0.2198865624746489
$ws.Range("E18").Value = 0.2694167480348568
$ws.Range("F18").Value = 1.269146211057006
$ws.Range("G18").Value = 0.002435784481438765
$ws.Range("I18").Value = 0.9148520523860384
$ws.Range("J18").Value = 0.3943653500969475
$ws.Range("L18").Value = 0.6833070323532979
$ws.Range("M18").Value = 0.5128371464040242
$ws.Range("O18").Value = 3.170748903436532
# Row 19
$ws.Range("B19").Value = 1.547062692437521
$ws.Range("D19").Value = 0.2198020778225853
$ws.Range("E19").Value = 0.268856204459901
$ws.Range("F19").Value = 1.269827867331593
$ws.Range("G19").Value = 0.002436041795017759
$ws.Range("I19").Value = 0.9160600081605317
$ws.Range("J19").Value = 0.3930987516870488
$ws.Range("L19").Value = 0.6789138263907262
$ws.Range("M19").Value = 0.5097671468665936
$ws.Range("O19").Value = 3.170456735724741
# Row 20
$ws.Range("B20").Value = 1.58640675843094
$ws.Range("D20").Value = 0.220184545867177
$ws.Range("E20").Value = 0.2713805592059657
$ws.Range("F20").Value = 1.266825147583141
$ws.Range("G20").Value = 0.002434891064048709
$ws.Range("I20").Value = 0.9106709897425951
$ws.Range("J20").Value = 0.3988001231138014
$ws.Range("L20").Value = 0.6986792754789803
$ws.Range("M20").Value = 0.5235798840071624
$ws.Range("O20").Value = 3.171941088835212
# Row 21
$ws.Range("B21").Value = 1.718483328938191
$ws.Range("D21").Value = 0.2215091371315125
$ws.Range("E21").Value = 0.2798934951658509
$ws.Range("F21").Value = 1.257872280461164
$ws.Range("G21").Value = 0.002431150600523198
$ws.Range("I21").Value = 0.8933862810834157
$ws.Range("J21").Value = 0.4179811921248699
$ws.Range("L21").Value = 0.7650066660186781
$ws.Range("M21").Value = 0.5699406906786777
$ws.Range("O21").Value = 3.179901346897282
# Row 22
$ws.Range("B22").Value = 1.804709311782801
$ws.Range("D22").Value = 0.2224039392984594
$ws.Range("E22").Value = 0.2854795947153974
$ws.Range("F22").Value = 1.252869680079996
$ws.Range("G22").Value = 0.002428799381151646
$ws.Range("I22").Value = 0.8827028588857253
$ws.Range("J22").Value = 0.4305339363081515
$ws.Range("L22").Value = 0.808290088768473
$ws.Range("M22").Value = 0.6002012350800214
$ws.Range("O22").Value = 3.187326471707223
# Row 23
$ws.Range("B23").Value = 1.758698006624229
$ws.Range("D23").Value = 0.2219237272541079
$ws.Range("E23").Value = 0.2824962045277744
$ws.Range("F23").Value = 1.255461882636709
$ws.Range("G23").Value = 0.002430045751277596
$ws.Range("I23").Value = 0.8883486521772177
$ws.Range("J23").Value = 0.423832890841922
$ws.Range("L23").Value = 0.7851951721447961
$ws.Range("M23").Value = 0.5840543707687829
$ws.Range("O23").Value = 3.183159735984162
# Row 24
$ws.Range("B24").Value = 1.584246577374131
$ws.Range("D24").Value = 0.2201633928041531
$ws.Range("E24").Value = 0.2712418127363421
$ws.Range("F24").Value = 1.266985774672321
$ws.Range("G24").Value = 0.002434953782407656
$ws.Range("I24").Value = 0.9109638396557393
$ws.Range("J24").Value = 0.3984869309346237
$ws.Range("L24").Value = 0.6975941479878998
$ws.Range("M24").Value = 0.522821527543357
$ws.Range("O24").Value = 3.17184840790955
# Row 25
$ws.Range("B25").Value = 1.395725865916859
$ws.Range("D25").Value = 0.2183919610284519
$ws.Range("E25").Value = 0.259205012320983
$ws.Range("F25").Value = 1.283046298340068
$ws.Range("G25").Value = 0.002440651435282986
$ws.Range("I25").Value = 0.9379856743496617
$ws.Range("J25").Value = 0.3712323244072309
$ws.Range("L25").Value = 0.6028498480156657
$ws.Range("M25").Value = 0.4566246594188
$ws.Range("O25").Value = 3.1691456204058
